$d = $word.ActiveDocument

function Set-CalibriFont($rng) {
    $rng.Font.NameAscii   = "Calibri"
    $rng.Font.NameFarEast = "Calibri"
    $rng.Font.NameOther   = "Calibri"
    $rng.Font.NameBi      = "Calibri"
}

# ---------------------------------------------------------------------
# Hunk 1: Append three new runs after "... for my admin panel."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("for my admin panel.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

$rng.InsertAfter(" As for ip whitelisting I have created an if statement within my ws.php file to check if the user’s domain url matches the parameters within the If statement and if they don’t then the request dies and they are unable to access it.")
Set-CalibriFont $rng

$rng.Collapse(0)
$rng.InsertAfter(" I have also got an if statement that checks the access rights of the user ")
Set-CalibriFont $rng

$rng.Collapse(0)
$rng.InsertAfter("when they try to login to the admin panel located within the adminLogin function.")
Set-CalibriFont $rng

# ---------------------------------------------------------------------
# Hunk 2: Split the "PHP." run into "PHP" + new runs + "."
# ---------------------------------------------------------------------
$rngPhp = $d.Content
$rngPhp.Find.Execute("PHP.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$phpStart = $rngPhp.Start
$phpEnd = $rngPhp.End

# Isolate "PHP" from "." into two separate runs via a formatting
# round-trip (forces a run split at both ends of the sub-range without
# touching the text itself, so the original run's identity is kept).
$rngPhpOnly = $d.Range($phpStart, $phpEnd - 1)
$rngPhpOnly.Bold = $true
$rngPhpOnly.Bold = $false

# Remove the now-isolated trailing "." run (it sits at the very end of
# the paragraph, so deleting it does not disturb any sibling runs).
$rngDot = $d.Range($phpEnd - 1, $phpEnd)
$rngDot.Delete()

# Re-append the replacement tail at the (now) end of the paragraph.
$rngTail = $d.Range($phpEnd - 1, $phpEnd - 1)
$rngTail.InsertAfter(" with custom code to check if empty and checkValidity")
Set-CalibriFont $rngTail

$rngTail.Collapse(0)
$rngTail.InsertAfter(".")
Set-CalibriFont $rngTail
